$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 3).Value = 8.511166290682612
$ws.Cells.Item(2, 4).Value = 7.336602992512772
$ws.Cells.Item(2, 5).Value = 12.17420243942703
$ws.Cells.Item(2, 6).Value = 39.4514277408683
$ws.Cells.Item(2, 7).Value = 3.684023309629047
$ws.Cells.Item(2, 9).Value = 30.53185555868647
$ws.Cells.Item(2, 10).Value = 9.850659160572398
$ws.Cells.Item(2, 12).Value = 9.765880897118548
$ws.Cells.Item(2, 13).Value = 26.67294745005681
$ws.Cells.Item(2, 14).Value = 17.61622787179599
$ws.Cells.Item(2, 15).Value = 30.77590289894916
$ws.Cells.Item(3, 3).Value = 8.539426696961097
$ws.Cells.Item(3, 4).Value = 7.35233068474836
$ws.Cells.Item(3, 5).Value = 12.22110154123809
$ws.Cells.Item(3, 6).Value = 39.39547983692009
$ws.Cells.Item(3, 7).Value = 3.687202568247016
$ws.Cells.Item(3, 9).Value = 30.51466231771913
$ws.Cells.Item(3, 10).Value = 9.882932442996141
$ws.Cells.Item(3, 12).Value = 9.794199200526787
$ws.Cells.Item(3, 13).Value = 26.06779840092776
$ws.Cells.Item(3, 14).Value = 17.34977537613874
$ws.Cells.Item(3, 15).Value = 30.73476350849494
$ws.Cells.Item(4, 3).Value = 8.557949298008847
$ws.Cells.Item(4, 4).Value = 7.362543191737735
$ws.Cells.Item(4, 5).Value = 12.251318481412
$ws.Cells.Item(4, 6).Value = 39.37167767804947
$ws.Cells.Item(4, 7).Value = 3.689258618311479
$ws.Cells.Item(4, 9).Value = 30.51177588678313
$ws.Cells.Item(4, 10).Value = 9.903699383792274
$ws.Cells.Item(4, 12).Value = 9.812428068074746
$ws.Cells.Item(4, 13).Value = 25.68958484395487
$ws.Cells.Item(4, 14).Value = 17.18612247847599
$ws.Cells.Item(4, 15).Value = 30.71741959624904
$ws.Cells.Item(5, 3).Value = 8.565791901223147
$ws.Cells.Item(5, 4).Value = 7.36684488146725
$ws.Cells.Item(5, 5).Value = 12.26399032842657
$ws.Cells.Item(5, 6).Value = 39.36463444889758
$ws.Cells.Item(5, 7).Value = 3.690122711189841
$ws.Cells.Item(5, 9).Value = 30.51252812696292
$ws.Cells.Item(5, 10).Value = 9.91240192504825
$ws.Cells.Item(5, 12).Value = 9.820068666399969
$ws.Cells.Item(5, 13).Value = 25.53397500618296
$ws.Cells.Item(5, 14).Value = 17.11949846769474
$ws.Cells.Item(5, 15).Value = 30.71234496560293
$ws.Cells.Item(6, 3).Value = 8.567111950329483
$ws.Cells.Item(6, 4).Value = 7.367567638598743
$ws.Cells.Item(6, 5).Value = 12.26611614126723
$ws.Cells.Item(6, 6).Value = 39.36362541565787
$ws.Cells.Item(6, 7).Value = 3.690267780408633
$ws.Cells.Item(6, 9).Value = 30.5127694767217
$ws.Cells.Item(6, 10).Value = 9.913861482076316
$ws.Cells.Item(6, 12).Value = 9.821350216036391
$ws.Cells.Item(6, 13).Value = 25.50805208898878
$ws.Cells.Item(6, 14).Value = 17.10844188927237
$ws.Cells.Item(6, 15).Value = 30.71162275216019
$ws.Cells.Item(7, 3).Value = 8.558053873517565
$ws.Cells.Item(7, 4).Value = 7.362600638541752
$ws.Cells.Item(7, 5).Value = 12.25148792678556
$ws.Cells.Item(7, 6).Value = 39.37157193254667
$ws.Cells.Item(7, 7).Value = 3.689270165422382
$ws.Cells.Item(7, 9).Value = 30.51177822486952
$ws.Cells.Item(7, 10).Value = 9.903815777263322
$ws.Cells.Item(7, 12).Value = 9.812530251890667
$ws.Cells.Item(7, 13).Value = 25.68749199279088
$ws.Cells.Item(7, 14).Value = 17.18522359277541
$ws.Cells.Item(7, 15).Value = 30.71734308568654
$ws.Cells.Item(8, 3).Value = 8.520667603047251
$ws.Cells.Item(8, 4).Value = 7.341910752474484
$ws.Cells.Item(8, 5).Value = 12.19007895852113
$ws.Cells.Item(8, 6).Value = 39.42995001330954
$ws.Cells.Item(8, 7).Value = 3.685097993427248
$ws.Cells.Item(8, 9).Value = 30.52433559619864
$ws.Cells.Item(8, 10).Value = 9.861590031061953
$ws.Cells.Item(8, 12).Value = 9.775470815216519
$ws.Cells.Item(8, 13).Value = 26.46578170889866
$ws.Cells.Item(8, 14).Value = 17.52441218245954
$ws.Cells.Item(8, 15).Value = 30.76007720497111
$ws.Cells.Item(9, 3).Value = 8.456635179770668
$ws.Cells.Item(9, 4).Value = 7.305733437102417
$ws.Cells.Item(9, 5).Value = 12.08088465536554
$ws.Cells.Item(9, 6).Value = 39.62788363070951
$ws.Cells.Item(9, 7).Value = 3.67773720313033
$ws.Cells.Item(9, 9).Value = 30.60977660153833
$ws.Cells.Item(9, 10).Value = 9.786298692085772
$ws.Cells.Item(9, 12).Value = 9.709443568389375
$ws.Cells.Item(9, 13).Value = 27.93149861200899
$ws.Cells.Item(9, 14).Value = 18.18577634442784
$ws.Cells.Item(9, 15).Value = 30.90649689737817
$ws.Cells.Item(10, 3).Value = 8.415240399112736
$ws.Cells.Item(10, 4).Value = 7.281814948574279
$ws.Cells.Item(10, 5).Value = 12.00744199765367
$ws.Cells.Item(10, 6).Value = 39.82368649589556
$ws.Cells.Item(10, 7).Value = 3.672823831914902
$ws.Cells.Item(10, 9).Value = 30.70948299392133
$ws.Cells.Item(10, 10).Value = 9.735515929659746
$ws.Cells.Item(10, 12).Value = 9.664943451618393
$ws.Cells.Item(10, 13).Value = 28.96182682874153
$ws.Cells.Item(10, 14).Value = 18.66495500474066
$ws.Cells.Item(10, 15).Value = 31.05188237897725
$ws.Cells.Item(11, 3).Value = 8.397633985258807
$ws.Cells.Item(11, 4).Value = 7.271507682478067
$ws.Cells.Item(11, 5).Value = 11.97549036417444
$ws.Cells.Item(11, 6).Value = 39.92354534172868
$ws.Cells.Item(11, 7).Value = 3.670694775409819
$ws.Cells.Item(11, 9).Value = 30.76280161656059
$ws.Cells.Item(11, 10).Value = 9.713388059351756
$ws.Cells.Item(11, 12).Value = 9.645561046891304
$ws.Cells.Item(11, 13).Value = 29.41861924795643
$ws.Cells.Item(11, 14).Value = 18.88059259776184
$ws.Cells.Item(11, 15).Value = 31.12611834762259
$ws.Cells.Item(12, 3).Value = 8.391142862748373
$ws.Cells.Item(12, 4).Value = 7.267686750245201
$ws.Cells.Item(12, 5).Value = 11.96359975369159
$ws.Cells.Item(12, 6).Value = 39.96289320156857
$ws.Cells.Item(12, 7).Value = 3.669903714108718
$ws.Cells.Item(12, 9).Value = 30.78412930265667
$ws.Cells.Item(12, 10).Value = 9.705148054164251
$ws.Cells.Item(12, 12).Value = 9.638344569749766
$ws.Cells.Item(12, 13).Value = 29.5897445717574
$ws.Cells.Item(12, 14).Value = 18.96184110352726
$ws.Cells.Item(12, 15).Value = 31.15538171914881
$ws.Cells.Item(13, 3).Value = 8.392533012493967
$ws.Cells.Item(13, 4).Value = 7.268506004106248
$ws.Cells.Item(13, 5).Value = 11.96615133732937
$ws.Cells.Item(13, 6).Value = 39.95435104097167
$ws.Cells.Item(13, 7).Value = 3.670073410144133
$ws.Cells.Item(13, 9).Value = 30.77948557150611
$ws.Cells.Item(13, 10).Value = 9.70691649977385
$ws.Cells.Item(13, 12).Value = 9.639893296235197
$ws.Cells.Item(13, 13).Value = 29.5529741471459
$ws.Cells.Item(13, 14).Value = 18.94436207452006
$ws.Cells.Item(13, 15).Value = 31.14902832391053
$ws.Cells.Item(14, 3).Value = 8.397096428511215
$ws.Cells.Item(14, 4).Value = 7.271191685947484
$ws.Cells.Item(14, 5).Value = 11.97450793708024
$ws.Cells.Item(14, 6).Value = 39.92675189041364
$ws.Cells.Item(14, 7).Value = 3.670629390853601
$ws.Cells.Item(14, 9).Value = 30.76453351918515
$ws.Cells.Item(14, 10).Value = 9.712707360275765
$ws.Cells.Item(14, 12).Value = 9.644964876797538
$ws.Cells.Item(14, 13).Value = 29.43273570947823
$ws.Cells.Item(14, 14).Value = 18.88728555387434
$ws.Cells.Item(14, 15).Value = 31.1285028583581
$ws.Cells.Item(15, 3).Value = 8.399914578749664
$ws.Cells.Item(15, 4).Value = 7.272847442033585
$ws.Cells.Item(15, 5).Value = 11.97965375989746
$ws.Cells.Item(15, 6).Value = 39.91004575247018
$ws.Cells.Item(15, 7).Value = 3.670971917656727
$ws.Cells.Item(15, 9).Value = 30.75552278057336
$ws.Cells.Item(15, 10).Value = 9.71627255656716
$ws.Cells.Item(15, 12).Value = 9.648087395892379
$ws.Cells.Item(15, 13).Value = 29.35884104885253
$ws.Cells.Item(15, 14).Value = 18.85226922456105
$ws.Cells.Item(15, 15).Value = 31.11608000818318
$ws.Cells.Item(16, 3).Value = 8.416415654602948
$ws.Cells.Item(16, 4).Value = 7.28250006940494
$ws.Cells.Item(16, 5).Value = 12.00955937384919
$ws.Cells.Item(16, 6).Value = 39.81737609244581
$ws.Cells.Item(16, 7).Value = 3.67296509768085
$ws.Cells.Item(16, 9).Value = 30.70615807254783
$ws.Cells.Item(16, 10).Value = 9.736981571460408
$ws.Cells.Item(16, 12).Value = 9.66622741158217
$ws.Cells.Item(16, 13).Value = 28.93172254927511
$ws.Cells.Item(16, 14).Value = 18.65080921979229
$ws.Cells.Item(16, 15).Value = 31.04719280421473
$ws.Cells.Item(17, 3).Value = 8.426852090507211
$ws.Cells.Item(17, 4).Value = 7.288568323250813
$ws.Cells.Item(17, 5).Value = 12.02827826943062
$ws.Cells.Item(17, 6).Value = 39.76327727709418
$ws.Cells.Item(17, 7).Value = 3.674214952959941
$ws.Cells.Item(17, 9).Value = 30.67790851576591
$ws.Cells.Item(17, 10).Value = 9.749934749124577
$ws.Cells.Item(17, 12).Value = 9.677575812081409
$ws.Cells.Item(17, 13).Value = 28.66654620384648
$ws.Cells.Item(17, 14).Value = 18.52656781121208
$ws.Cells.Item(17, 15).Value = 31.00699858178627
$ws.Cells.Item(18, 3).Value = 8.432970087153564
$ws.Cells.Item(18, 4).Value = 7.292112605490553
$ws.Cells.Item(18, 5).Value = 12.03918214846853
$ws.Cells.Item(18, 6).Value = 39.73317762021221
$ws.Cells.Item(18, 7).Value = 3.674943824657504
$ws.Cells.Item(18, 9).Value = 30.66241006963048
$ws.Cells.Item(18, 10).Value = 9.757476739504199
$ws.Cells.Item(18, 12).Value = 9.684184180423415
$ws.Cells.Item(18, 13).Value = 28.51291048577084
$ws.Cells.Item(18, 14).Value = 18.45489099538647
$ws.Cells.Item(18, 15).Value = 30.98464316154891
$ws.Cells.Item(19, 3).Value = 8.435061327691914
$ws.Cells.Item(19, 4).Value = 7.293321917725909
$ws.Cells.Item(19, 5).Value = 12.04289761603728
$ws.Cells.Item(19, 6).Value = 39.72316148682224
$ws.Cells.Item(19, 7).Value = 3.675192326143345
$ws.Cells.Item(19, 9).Value = 30.65729158065073
$ws.Cells.Item(19, 10).Value = 9.760046091114374
$ws.Cells.Item(19, 12).Value = 9.686435601731176
$ws.Cells.Item(19, 13).Value = 28.46070529389354
$ws.Cells.Item(19, 14).Value = 18.43058752159059
$ws.Cells.Item(19, 15).Value = 30.97720544987973
$ws.Cells.Item(20, 3).Value = 8.425729187813415
$ws.Cells.Item(20, 4).Value = 7.287916762005524
$ws.Cells.Item(20, 5).Value = 12.02627141081284
$ws.Cells.Item(20, 6).Value = 39.76893110004381
$ws.Cells.Item(20, 7).Value = 3.67408087063343
$ws.Cells.Item(20, 9).Value = 30.68083816423432
$ws.Cells.Item(20, 10).Value = 9.74854637895864
$ws.Cells.Item(20, 12).Value = 9.676359369433378
$ws.Cells.Item(20, 13).Value = 28.69489094643336
$ws.Cells.Item(20, 14).Value = 18.53981645249564
$ws.Cells.Item(20, 15).Value = 31.01119842152963
$ws.Cells.Item(21, 3).Value = 8.395751264813921
$ws.Cells.Item(21, 4).Value = 7.270400606982208
$ws.Cells.Item(21, 5).Value = 11.97204774168186
$ws.Cells.Item(21, 6).Value = 39.93481696452243
$ws.Cells.Item(21, 7).Value = 3.670465674912256
$ws.Cells.Item(21, 9).Value = 30.76889450870923
$ws.Cells.Item(21, 10).Value = 9.711002667659782
$ws.Cells.Item(21, 12).Value = 9.643471891134006
$ws.Cells.Item(21, 13).Value = 29.46810397504728
$ws.Cells.Item(21, 14).Value = 18.90406196233941
$ws.Cells.Item(21, 15).Value = 31.13450053602269
$ws.Cells.Item(22, 3).Value = 8.377185084260288
$ws.Cells.Item(22, 4).Value = 7.259431834030335
$ws.Cells.Item(22, 5).Value = 11.93782597816759
$ws.Cells.Item(22, 6).Value = 40.05216139936503
$ws.Cells.Item(22, 7).Value = 3.668191296520204
$ws.Cells.Item(22, 9).Value = 30.83306755888663
$ws.Cells.Item(22, 10).Value = 9.68727756103968
$ws.Cells.Item(22, 12).Value = 9.62269596985449
$ws.Cells.Item(22, 13).Value = 29.96259118160036
$ws.Cells.Item(22, 14).Value = 19.13970376443678
$ws.Cells.Item(22, 15).Value = 31.22179212749069
$ws.Cells.Item(23, 3).Value = 8.387000313101055
$ws.Cells.Item(23, 4).Value = 7.26524232169648
$ws.Cells.Item(23, 5).Value = 11.9559797428872
$ws.Cells.Item(23, 6).Value = 39.98872197236327
$ws.Cells.Item(23, 7).Value = 3.669397118053854
$ws.Cells.Item(23, 9).Value = 30.79821412791334
$ws.Cells.Item(23, 10).Value = 9.699866021546839
$ws.Cells.Item(23, 12).Value = 9.633718968794685
$ws.Cells.Item(23, 13).Value = 29.6997103056505
$ws.Cells.Item(23, 14).Value = 19.01418106339046
$ws.Cells.Item(23, 15).Value = 31.17459400342481
$ws.Cells.Item(24, 3).Value = 8.426236484726193
$ws.Cells.Item(24, 4).Value = 7.288211159762889
$ws.Cells.Item(24, 5).Value = 12.02717826878842
$ws.Cells.Item(24, 6).Value = 39.76637188275451
$ws.Cells.Item(24, 7).Value = 3.674141457122015
$ws.Cells.Item(24, 9).Value = 30.67951135644785
$ws.Cells.Item(24, 10).Value = 9.749173765030877
$ws.Cells.Item(24, 12).Value = 9.676909061388587
$ws.Cells.Item(24, 13).Value = 28.68207996017938
$ws.Cells.Item(24, 14).Value = 18.53382751115729
$ws.Cells.Item(24, 15).Value = 31.00929732819241
$ws.Cells.Item(25, 3).Value = 8.472964884259101
$ws.Cells.Item(25, 4).Value = 7.315051730113212
$ws.Cells.Item(25, 5).Value = 12.10922903918157
$ws.Cells.Item(25, 6).Value = 39.56544421422817
$ws.Cells.Item(25, 7).Value = 3.679641217975443
$ws.Cells.Item(25, 9).Value = 30.58016071141443
$ws.Cells.Item(25, 10).Value = 9.805867472363802
$ws.Cells.Item(25, 12).Value = 9.726598479960948
$ws.Cells.Item(25, 13).Value = 27.54249863896129
$ws.Cells.Item(25, 14).Value = 18.00772915515927
$ws.Cells.Item(25, 15).Value = 30.86021316994051
